$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(381, 20081400, 1641133000000, 1635886000000, -5247000000),
    @(382, 20081500, 0, 0, 0),
    @(383, 20081600, 0, 0, 0),
    @(384, 20081700, 1635886000000, 1663881000000, 27995000000),
    @(385, 20081800, 1663881000000, 1655467000000, -8414000000),
    @(386, 20081900, 1655467000000, 1636393000000, -19074000000)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $ws.Cells.Item($row, 1).Value = $entry[1]
    $ws.Cells.Item($row, 2).Value = $entry[2]
    $ws.Cells.Item($row, 3).Value = $entry[3]
    $ws.Cells.Item($row, 4).Value = $entry[4]
}
